$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.302010295839347
$ws.Range("B3").Value = 0.16849896177434848
$ws.Range("B4").Value = 72.56810916183356

$ws.Range("A5").ClearContents()
$ws.Range("B5").ClearContents()
